$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: first_median/second_median and open/close columns updated to 1288
$ws.Range("C10").Value = 1288
$ws.Range("E10").Value = 1288

# Row 31: first_median and open columns updated to 719
$ws.Range("B31").Value = 719
$ws.Range("D31").Value = 719
